$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "0.1.5" version history row (row 8) ---
$ws.Range("A8").Value = "0.1.5"
$ws.Range("B8").Value = "AUTOMATA CELULAR - copia (12)"
$ws.Range("C8").Value = "-Agrupation and desagrupation to be done in functions.`n-Change reproduction and distribution to two parts.`n-UI: Delete rows according to working functionality.`n-UI: condicionate IF to be associated.`n-Implement mutations.`n-With 4 or less niches the distribution is not equaly done.`n-Document every function.`n*When an actor dies, the recipient should deassociate.`n-UI: automatically fill aggrupation data.`n*Agrupation and association to itself do not work properly."
$ws.Range("D8").Value = "-SI performace improvement.`n-Save and load implemented to file temp.csv.`n-Consumption can be a float due to the percentage of DF to consume.`n-The program can work now with 1 niches or more.`n-Relative individual selection pressure added."
$ws.Range("E8").Value = "Python 3.6.1"
$ws.Range("F8").Value = "Qt version: 5.6.2`nSIP version: 4.18`nPyQt version: 5.6"
$ws.Range("G8").Value = " PyInstaller 3.3.1"

# Row grew to fit the ten wrapped lines of the change log (10 x 14.4pt)
$ws.Rows.Item(8).RowHeight = 144

# --- Update sheet view: keep header row frozen, move the active selection ---
$win = $wb.Windows.Item(1)
[void]($win.FreezePanes = $true)

[void]$ws.Range("D9").Select()
